$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Append " (Adjust.com, 2024)" right after the "Donations from the public
#    ..." bullet (paragraph about donations).
# ---------------------------------------------------------------------------
$donationsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Donations from the public or users who want to fund development towards bettering the service*") {
        $donationsPara = $p
        break
    }
}
$donationsPara.Range.InsertAfter(" (Adjust.com, 2024)")

# ---------------------------------------------------------------------------
# 2) Append " (AppsFlyer, 2024)" right after the "Cost per view ..." bullet.
# ---------------------------------------------------------------------------
$cpvPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Cost per view*") {
        $cpvPara = $p
        break
    }
}
$cpvPara.Range.InsertAfter(" (AppsFlyer, 2024)")

# ---------------------------------------------------------------------------
# 3) Add a References section at the end of the document.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertAfter("`r`rReferences:`rAdjust.com. (2024). A beginner" + [char]0x2019 + "s guide to app monetization | Adjust. [online] Available at: https://www.adjust.com/blog/how-to-monetize-your-app/?utm_source=chatgpt.com [Accessed 26 Jan. 2025].`rAppsFlyer (2024). In-app advertising done right " + [char]0x2013 + " the complete guide. [online] AppsFlyer. Available at: https://www.appsflyer.com/resources/guides/in-app-advertising/?utm_source=chatgpt.com [Accessed 26 Jan. 2025].")

# Re-fetch the paragraph list: the new paragraphs are the last 4 of the document.
$count = $d.Paragraphs.Count
$adjustPara = $d.Paragraphs.Item($count - 1)
$appsflyerPara = $d.Paragraphs.Item($count)

# --- Italicise the Adjust reference title ---
$adjustText = $adjustPara.Range.Text
$titleAdjust = "A beginner" + [char]0x2019 + "s guide to app monetization | Adjust"
$idxAdjust = $adjustText.IndexOf($titleAdjust)
$startAdjust = $adjustPara.Range.Start + $idxAdjust
$endAdjust = $startAdjust + $titleAdjust.Length
$italicRangeAdjust = $d.Range($startAdjust, $endAdjust)
$italicRangeAdjust.Font.Italic = $true

# --- Italicise the AppsFlyer reference title ---
$appsflyerText = $appsflyerPara.Range.Text
$titleAppsFlyer = "In-app advertising done right " + [char]0x2013 + " the complete guide"
$idxAppsFlyer = $appsflyerText.IndexOf($titleAppsFlyer)
$startAppsFlyer = $appsflyerPara.Range.Start + $idxAppsFlyer
$endAppsFlyer = $startAppsFlyer + $titleAppsFlyer.Length
$italicRangeAppsFlyer = $d.Range($startAppsFlyer, $endAppsFlyer)
$italicRangeAppsFlyer.Font.Italic = $true

Write-Output "done"
